$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add E1 "MSE" and F1 "MAE" with the same header formatting
# (bold/border/centered) already used by A1:D1. Copy format from an existing
# header cell (PasteSpecial formats-only) then set the text.
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1").Value = "MSE"
$ws.Range("F1").Value = "MAE"

# --- New data rows 3-8: repeated "number_of_seasons" rows with MSE/MAE only
# (Precisao/Desvio Padrao columns C/D are left blank for these rows).
for ($r = 3; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = 100
    $ws.Cells.Item($r, 2).Value = "number_of_seasons"
    $ws.Cells.Item($r, 5).Value = 0.2213313391999655
    $ws.Cells.Item($r, 6).Value = 0.4597599784197205
}
